$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.421.47'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.26%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = '3.488.93'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = '''553.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.23%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = '''178.28'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.25%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = '''0.636'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.80%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E9').Value = '  -1.20%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('E10').Value = '  +1.69%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = '''53.51'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.78%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').Value = '  -1.94%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').Value = '  -2.72%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = '4.043.58'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.12%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('B15').Style = 'Normal'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('C15').Style = 'Normal'
$ws.Range('D15').Value = '3.486.15'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.19%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('B16').Value = 'TRON'
$ws.Range('B16').Style = 'Normal'
$ws.Range('C16').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('C16').Style = 'Normal'
$ws.Range('D16').Value = '''0.121'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('B17').Style = 'Normal'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C17').Style = 'Normal'
$ws.Range('D17').Value = '''18.44'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.69%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = '''12.03'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.63%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('D19').Value = '65.422.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.95%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = '''0.989'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.10%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = '''417.01'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.12%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = '  +1.81%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = '''86.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.85%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = '''4.11'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.43%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = '''12.71'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.09%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = '  -11.03%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').Value = '  -3.31%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = '''6.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.03%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = '''8.98'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.80%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = '''30.22'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.50%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').Value = '  -5.97%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = '''608.75'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -11.38%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = '''11.68'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = '  +9.90%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = '''0.999'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = '''37.27'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.38%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('B39').Value = 'PEPE'
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = '0.0₃0786'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -6.02%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('B40').Value = 'Maker'
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = '3.375.94'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.42%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = '''0.380'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.98%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = '''0.999'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = '''3.24'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.55%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = '''2.84'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.74%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').Value = '  -9.96%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('B46').Value = 'ApeXProtocol'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = '''3.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.18%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('B47').Value = 'VeChain'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = '''0.0412'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.39%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = '''2.71'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.20%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = '  +0.99%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = '''8.47'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.02%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = '''137.78'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.23%  '
$ws.Range('E51').Style = 'Normal'
